$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.113956600348786
$ws.Range("C2").Value = 0.07041660960565821
$ws.Range("D2").Value = 0.002230355564519115
$ws.Range("E2").Value = 0.06147272907016665
$ws.Range("F2").Value = 4.984343643786957
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 3.280589368861399
$ws.Range("J2").Value = 0.1775483532455873
$ws.Range("K2").Value = 0.9994063224198442
$ws.Range("L2").Value = 0.326208326667448
$ws.Range("M2").Value = 0.3088152948878466

# Row 3
$ws.Range("B3").Value = 1.101395631268787
$ws.Range("C3").Value = 0.0653133134775743
$ws.Range("D3").Value = 0.002203931866427311
$ws.Range("E3").Value = 0.06170422803215736
$ws.Range("F3").Value = 4.944654813428286
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 3.255442672115777
$ws.Range("J3").Value = 0.1775156476803286
$ws.Range("K3").Value = 0.9807796220937064
$ws.Range("L3").Value = 0.3260158177419541
$ws.Range("M3").Value = 0.3069278121838401

# Row 4
$ws.Range("B4").Value = 1.094444579545865
$ws.Range("C4").Value = 0.06222300632200017
$ws.Range("D4").Value = 0.002190876950550091
$ws.Range("E4").Value = 0.06185837856750354
$ws.Range("F4").Value = 4.921493858767292
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 3.240666928315861
$ws.Range("J4").Value = 0.1775099148167723
$ws.Range("K4").Value = 0.9700475144178995
$ws.Range("L4").Value = 0.3260330249506893
$ws.Range("M4").Value = 0.3059494516714381

# Row 5
$ws.Range("B5").Value = 1.09180370621101
$ws.Range("C5").Value = 0.06097446241061277
$ws.Range("D5").Value = 0.002186354025598192
$ws.Range("E5").Value = 0.06192422509743567
$ws.Range("F5").Value = 4.912359305783085
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 3.2348126497252
$ws.Range("J5").Value = 0.1775111897703674
$ws.Range("K5").Value = 0.9658514350862788
$ws.Range("L5").Value = 0.3260741512716123
$ws.Range("M5").Value = 0.3055962387205007

# Row 6
$ws.Range("B6").Value = 1.091376779643184
$ws.Range("C6").Value = 0.06076779190237858
$ws.Range("D6").Value = 0.002185651140839795
$ws.Range("E6").Value = 0.0619353420725206
$ws.Range("F6").Value = 4.910860862385633
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 3.233850631598116
$ws.Range("J6").Value = 0.1775116196712645
$ws.Range("K6").Value = 0.9651653941052416
$ws.Range("L6").Value = 0.3260830425985759
$ws.Range("M6").Value = 0.3055403368529319

# Row 7
$ws.Range("B7").Value = 1.09440818711272
$ws.Range("C7").Value = 0.06220612445581253
$ws.Range("D7").Value = 0.002190812725126889
$ws.Range("E7").Value = 0.06185925432043327
$ws.Range("F7").Value = 4.921369437341355
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 3.240587299636744
$ws.Range("J7").Value = 0.1775099173856614
$ws.Range("K7").Value = 0.9699902064082266
$ws.Range("L7").Value = 0.3260334413809147
$ws.Range("M7").Value = 0.3059445038860638

# Row 8
$ws.Range("B8").Value = 1.109467693689879
$ws.Range("C8").Value = 0.06864799846454162
$ws.Range("D8").Value = 0.00222058690343907
$ws.Range("E8").Value = 0.06155006284918274
$ws.Range("F8").Value = 4.970408120038016
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 3.271780802485821
$ws.Range("J8").Value = 0.1775341001085238
$ws.Range("K8").Value = 0.9928376259279048
$ws.Range("L8").Value = 0.326113883363405
$ws.Range("M8").Value = 0.308127057352781

# Row 9
$ws.Range("B9").Value = 1.14503288877799
$ws.Range("C9").Value = 0.08162646807340934
$ws.Range("D9").Value = 0.002304129704372215
$ws.Range("E9").Value = 0.06103862800943816
$ws.Range("F9").Value = 5.076170379350856
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 3.338233256549344
$ws.Range("J9").Value = 0.1776952742083253
$ws.Range("K9").Value = 1.043232231780735
$ws.Range("L9").Value = 0.3273435713869475
$ws.Range("M9").Value = 0.3138372194715835

# Row 10
$ws.Range("B10").Value = 1.174835432528567
$ws.Range("C10").Value = 0.09137883456013185
$ws.Range("D10").Value = 0.002380874737957228
$ws.Range("E10").Value = 0.06072018855940353
$ws.Range("F10").Value = 5.159752352218817
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 3.390296884204432
$ws.Range("J10").Value = 0.177882962879778
$ws.Range("K10").Value = 1.083670504081368
$ws.Range("L10").Value = 0.3288977550276044
$ws.Range("M10").Value = 0.3189020754153837

# Row 11
$ws.Range("B11").Value = 1.189190207051553
$ws.Range("C11").Value = 0.09586397661198021
$ws.Range("D11").Value = 0.002419133039650134
$ws.Range("E11").Value = 0.06058765482392836
$ws.Range("F11").Value = 5.199059108563802
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 3.414690671700512
$ws.Range("J11").Value = 0.1779833785011249
$ws.Range("K11").Value = 1.10280968935092
$ws.Range("L11").Value = 0.3297455607443283
$ws.Range("M11").Value = 0.3213946448100131

# Row 12
$ws.Range("B12").Value = 1.194740505880702
$ws.Range("C12").Value = 0.09756948384833208
$ws.Range("D12").Value = 0.002434102115588033
$ws.Range("E12").Value = 0.06053923125741179
$ws.Range("F12").Value = 5.214128627864312
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 3.424030286818422
$ws.Range("J12").Value = 0.1780235633236806
$ws.Range("K12").Value = 1.110164144121001
$ws.Range("L12").Value = 0.3300867968577492
$ws.Range("M12").Value = 0.3223655774541641

# Row 13
$ws.Range("B13").Value = 1.193540062554945
$ws.Range("C13").Value = 0.09720185624911437
$ws.Range("D13").Value = 0.002430856839691131
$ws.Range("E13").Value = 0.0605495818133992
$ws.Range("F13").Value = 5.2108749113832
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 3.422014285028993
$ws.Range("J13").Value = 0.1780148127930481
$ws.Range("K13").Value = 1.108575479084294
$ws.Range("L13").Value = 0.3300124084666862
$ws.Range("M13").Value = 0.3221552676104409

# Row 14
$ws.Range("B14").Value = 1.189644540790425
$ws.Range("C14").Value = 0.09600414738289942
$ws.Range("D14").Value = 0.002420354903154731
$ws.Range("E14").Value = 0.06058363568085312
$ws.Range("F14").Value = 5.200295180442936
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 3.415456998001204
$ws.Range("J14").Value = 0.1779866412629048
$ws.Range("K14").Value = 1.103412603501425
$ws.Range("L14").Value = 0.3297732301175671
$ws.Range("M14").Value = 0.321473982232618

# Row 15
$ws.Range("B15").Value = 1.187273320676866
$ws.Range("C15").Value = 0.09527144046924718
$ws.Range("D15").Value = 0.002413984877323472
$ws.Range("E15").Value = 0.06060472414620399
$ws.Range("F15").Value = 5.19383887381116
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 3.411453786569922
$ws.Range("J15").Value = 0.177969666560351
$ws.Range("K15").Value = 1.100264107582859
$ws.Range("L15").Value = 0.3296293542119599
$ws.Range("M15").Value = 0.3210601962729385

# Row 16
$ws.Range("B16").Value = 1.173913342439391
$ws.Range("C16").Value = 0.09108670694729426
$ws.Range("D16").Value = 0.002378441830147793
$ws.Range("E16").Value = 0.0607290972629535
$ws.Range("F16").Value = 5.157209430303652
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 3.388716988366127
$ws.Range("J16").Value = 0.1778767028198551
$ws.Range("K16").Value = 1.082434670575566
$ws.Range("L16").Value = 0.3288451773459968
$ws.Range("M16").Value = 0.3187429688033276

# Row 17
$ws.Range("B17").Value = 1.165921519225833
$ws.Range("C17").Value = 0.08853204357777145
$ws.Range("D17").Value = 0.002357494634466661
$ws.Range("E17").Value = 0.06080854713337303
$ws.Range("F17").Value = 5.135067642617457
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 3.374950598697779
$ws.Range("J17").Value = 0.1778235216978317
$ws.Range("K17").Value = 1.071687293696812
$ws.Range("L17").Value = 0.3284001320819172
$ws.Range("M17").Value = 0.3173696706694074

# Row 18
$ws.Range("B18").Value = 1.161399905893006
$ws.Range("C18").Value = 0.08706725220189071
$ws.Range("D18").Value = 0.002345761342585462
$ws.Range("E18").Value = 0.06085540527226652
$ws.Range("F18").Value = 5.122453204132142
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 3.367099334348964
$ws.Range("J18").Value = 0.1777943486423048
$ws.Range("K18").Value = 1.065575682862885
$ws.Range("L18").Value = 0.3281574098761482
$ws.Range("M18").Value = 0.3165975358548252

# Row 19
$ws.Range("B19").Value = 1.159881865651357
$ws.Range("C19").Value = 0.0865720840266988
$ws.Range("D19").Value = 0.002341842741806488
$ws.Range("E19").Value = 0.06087147023627892
$ws.Range("F19").Value = 5.118202938131247
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 3.364452499339166
$ws.Range("J19").Value = 0.1777847143030939
$ws.Range("K19").Value = 1.063518420578418
$ws.Range("L19").Value = 0.3280775067788753
$ws.Range("M19").Value = 0.3163391547339742

# Row 20
$ws.Range("B20").Value = 1.166764494353231
$ws.Range("C20").Value = 0.08880351696635103
$ws.Range("D20").Value = 0.002359691898762506
$ws.Range("E20").Value = 0.0607999694946022
$ws.Range("F20").Value = 5.137412156690772
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 3.376409139839708
$ws.Range("J20").Value = 0.1778290364556732
$ws.Range("K20").Value = 1.072824126505623
$ws.Range("L20").Value = 0.3284461363685836
$ws.Range("M20").Value = 0.3175140239948497

# Row 21
$ws.Range("B21").Value = 1.190785645304402
$ws.Range("C21").Value = 0.09635575074526059
$ws.Range("D21").Value = 0.002423426507021986
$ws.Range("E21").Value = 0.06057358542409208
$ws.Range("F21").Value = 5.203397684530415
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 3.417380256792768
$ws.Range("J21").Value = 0.1779948573321413
$ws.Range("K21").Value = 1.104926165566752
$ws.Range("L21").Value = 0.329842935077032
$ws.Range("M21").Value = 0.3216733585956035

# Row 22
$ws.Range("B22").Value = 1.207151905424865
$ws.Range("C22").Value = 0.1013328813222643
$ws.Range("D22").Value = 0.002467887378497124
$ws.Range("E22").Value = 0.06043590893844897
$ws.Range("F22").Value = 5.247600868285474
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 3.444753121512463
$ws.Range("J22").Value = 0.1781158175844979
$ws.Range("K22").Value = 1.126529541851937
$ws.Range("L22").Value = 0.3308734838418417
$ws.Range("M22").Value = 0.3245493689795111

# Row 23
$ws.Range("B23").Value = 1.198355958633528
$ws.Range("C23").Value = 0.0986726891812566
$ws.Range("D23").Value = 0.002443900875606531
$ws.Range("E23").Value = 0.06050845172590513
$ws.Range("F23").Value = 5.223910129757797
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 3.430089135840888
$ws.Range("J23").Value = 0.1780501079316821
$ws.Range("K23").Value = 1.114942445047348
$ws.Range("L23").Value = 0.3303127130680608
$ws.Range("M23").Value = 0.3229999834913571

# Row 24
$ws.Range("B24").Value = 1.166383157914197
$ws.Range("C24").Value = 0.08868077164837018
$ws.Range("D24").Value = 0.002358697551549938
$ws.Range("E24").Value = 0.06080384376486947
$ws.Range("F24").Value = 5.13635184310715
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 3.375749536462195
$ws.Range("J24").Value = 0.1778265388683771
$ws.Range("K24").Value = 1.072309955066942
$ws.Range("L24").Value = 0.3284252968959862
$ws.Range("M24").Value = 0.3174487076749841

# Row 25
$ws.Range("B25").Value = 1.134766212050437
$ws.Range("C25").Value = 0.07807777416979889
$ws.Range("D25").Value = 0.002278833049633278
$ws.Range("E25").Value = 0.0611668817791795
$ws.Range("F25").Value = 5.046528571265952
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 3.319688466478851
$ws.Range("J25").Value = 0.1776395019327595
$ws.Range("K25").Value = 1.029000119570185
$ws.Range("L25").Value = 0.3268964371138097
$ws.Range("M25").Value = 0.3121395989906546

